$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.465.30"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "3.670.41"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "646.93"
$ws.Range("E5").Value = "  -4.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.36"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.145"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.10"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.439"
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000231"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "4.285.02"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.50"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").Value = "3.682.69"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").Value = "69.474.76"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.96"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.42"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "465.59"
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.75"
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("E22").Value = "  -1.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.46"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").Value = "3.815.02"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.73"
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.91"
$ws.Range("E28").Value = "  -2.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.60"
$ws.Range("E29").Value = "  -3.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.66"
$ws.Range("E30").Value = "  -4.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.99"
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.58"
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.44"
$ws.Range("E34").Value = "  -2.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.163"
$ws.Range("E35").Value = "  +3.42%  "
$ws.Range("D36").Value = "3.659.61"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.35"
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("E39").Value = "  -5.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "178.57"
$ws.Range("E40").Value = "  +4.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  -4.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0890"
$ws.Range("E43").Value = "  -1.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.925"
$ws.Range("E44").Value = "  -1.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.66"
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.71"
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.26"
$ws.Range("E47").Value = "  -3.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.97"
$ws.Range("E48").Value = "  -5.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.80"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000265"
$ws.Range("E50").Value = "  -3.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.05"
$ws.Range("E51").Value = "  -5.71%  "
